# hon_planning_and_dev.xlsx edit
#
# Summary of the change (from the commit diff):
#   1. The worksheet "פרופיל הון תכנון ופיתוח" (sheetId 7, 3rd tab) is
#      renamed to "משאב הון תכנון ופיתוח". Excel automatically rewrites
#      every formula on other sheets that references the old sheet name
#      (e.g. on "צללית זמנית").
#   2. That renamed sheet becomes the active tab of the workbook
#      (activeTab points at it / it gets tabSelected="1"), with the
#      active-cell selection on it set to B3.
#   3. The previously-active first sheet "מדד על תכנון" is no longer the
#      selected tab, and its own lingering selection moves to G9.
#   4. On "מדד על תכנון", the per-row formulas in columns R and S
#      (R2:R13 = N*10/$N$14, S2:S13 = O*10/$O$14) and the row of MAX()
#      summaries in row 14 (I14:T14) are re-entered as single range
#      fills, which makes Excel store them as shared formulas.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the sheet (sheetId 7 / 3rd tab in the workbook) ---
$wsHon = $wb.Worksheets.Item("פרופיל הון תכנון ופיתוח")
$wsHon.Name = "משאב הון תכנון ופיתוח"

# --- 2. Re-enter the formulas on "מדד על תכנון" as range fills so they ---
#        collapse into shared formulas, matching the target workbook.
$wsIndex = $wb.Worksheets.Item("מדד על תכנון")
$wsIndex.Range("R2:R13").Formula = "=N2*10/`$N`$14"
$wsIndex.Range("S2:S13").Formula = "=O2*10/`$O`$14"
$wsIndex.Range("I14:T14").Formula = "=MAX(I2:I13)"

# --- 3. Update the view state: "מדד על תכנון" keeps selection G9 but is ---
#        no longer the active tab.
$wsIndex.Activate()
$wsIndex.Range("G9").Select()

# --- 4. "משאב הון תכנון ופיתוח" becomes the active tab, selection B3. ---
$wsHon.Activate()
$wsHon.Range("B3").Select()
